# Add_PriceAgreementItem_FromComparisonScreen
# Adds a new "AddPriceAgrmnt_CompareScreen" test sheet (copy-like of the
# FavFolder sheet, plus a new "searchItem"/"UnitofMeasure" pair of columns),
# and nudges a couple of pre-existing sheets' selections, matching the
# state Excel leaves behind after interactively adding + using the sheet.

$wb = $excel.ActiveWorkbook

# --- Sheet5 (AddNonPriceAgr_GlobalCatalog): selection moved to F8 ---
$wsGlobalCatalog = $wb.Worksheets.Item("AddNonPriceAgr_GlobalCatalog")
$wsGlobalCatalog.Activate() | Out-Null
$wsGlobalCatalog.Range("F8").Select() | Out-Null

# --- Sheet6 (AddPriceAgrmnt_FavFolder): selection moved to D6 ---
$wsFavFolder = $wb.Worksheets.Item("AddPriceAgrmnt_FavFolder")
$wsFavFolder.Activate() | Out-Null
$wsFavFolder.Range("D6").Select() | Out-Null

# --- New sheet7 (AddPriceAgrmnt_CompareScreen), appended as the last tab ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsCompare = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsCompare.Name = "AddPriceAgrmnt_CompareScreen"

$wsCompare.Range("A1").Value = "Role"
$wsCompare.Range("B1").Value = "Location          "
$wsCompare.Range("C1").Value = "searchItem"
$wsCompare.Range("D1").Value = "UnitofMeasure"

$wsCompare.Range("A2").Value = "REQUESTOR"
$wsCompare.Range("B2").Value = "XEEVA -MJ"
$wsCompare.Range("C2").Value = "Iphones"
$wsCompare.Range("D2").Value = "CU-CUBIC"

$wsCompare.Columns.Item(1).ColumnWidth = 11.42578125
$wsCompare.Columns.Item(2).ColumnWidth = 10.140625
$wsCompare.Columns.Item(3).ColumnWidth = 9.85546875
$wsCompare.Columns.Item(4).ColumnWidth = 14.42578125

# New sheet becomes the active / tab-selected sheet, cursor parked at D6,
# same as the other freshly-recorded test sheets.
$wsCompare.Activate() | Out-Null
$wsCompare.Range("D6").Select() | Out-Null
